$wb = $excel.ActiveWorkbook

$metaWs = $wb.Worksheets.Item("Metadata")
$elemWs = $wb.Worksheets.Item("Elements")

# --- Metadata sheet updates ---
# Version: 1.8.7 -> 1.8.11 (force text so Excel doesn't auto-convert "1.8.11" to a date serial)
$metaWs.Range("B3").NumberFormat = "@"
$metaWs.Range("B3").Value = "1.8.11"
# Date: 2024-05-08T11:46:20-04:00 -> 2024-06-13T17:23:26-04:00
$metaWs.Range("B8").Value = "2024-06-13T17:23:26-04:00"

# --- Elements sheet updates ---
# Row 6 is the "Extension.value[x]" slicing-definition row.
# Type(s) column (K) now lists the full set of allowed data types.
$elemWs.Range("K6").Value = "base64Binary`nbooleancanonicalcodedatedateTimedecimalidinstantintegermarkdownoidpositiveIntstringtimeunsignedInturiurluuidAddressAgeAnnotationAttachmentCodeableConceptCodingContactPointCountDistanceDurationHumanNameIdentifierMoneyPeriodQuantityRangeRatioReferenceSampledDataSignatureTimingContactDetailContributorDataRequirementExpressionParameterDefinitionRelatedArtifactTriggerDefinitionUsageContextDosageMeta"

# Slicing Rules column (AE) changes from "closed" to "open" since a value type slice was added.
$elemWs.Range("AE6").Value = "open"

# Column K needs to widen to fit the long Type(s) text (results in OOXML width="255").
$elemWs.Columns.Item(11).ColumnWidth = 254.17
